$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.297769546508789
$ws.Range("B1").Value = 2.280836343765259
$ws.Range("C1").Value = 2.882692098617554
$ws.Range("D1").Value = 3.319928884506226
$ws.Range("E1").Value = 1.815756320953369
